# Rename recipe in the sheet: "Rice and X" -> "Lomo Saltado with X", and add three new
# recipe columns (chickpeas moves from D, then mushrooms/salmon/tofu are appended as E/F/G)
# with their corresponding impact figures for the final delivery dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): update existing labels, add new ones for E1:G1 ---
$headers = @(
  'Recipe Lomo Saltado with beef',
  'Recipe Lomo Saltado with chicken',
  'Recipe Lomo Saltado with chickpeas',
  'Recipe Lomo Saltado with mushrooms',
  'Recipe Lomo Saltado with salmon',
  'Recipe Lomo Saltado with tofu'
)

# The three brand-new header cells (E1:G1) need the same bold/centered header style
# already used by A1:D1 - copy the format from D1 onto them before writing text.
$ws.Range("D1").Copy()
$ws.Range("E1:G1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

for ($i = 0; $i -lt $headers.Length; $i++) {
  $ws.Cells.Item(1, $i + 2).Value = $headers[$i]
}

# --- Data rows 2-15: columns B..G (B..D already existed, E..G are brand new) ---
$data = @{
  2 = @(13.7190236, 6.9636386, 8.8705436, 4.4898036, 5.6031886, 5.513553600000001);
  3 = @(3.0133508, 1.5405158, 2.4227608, 0.7474508, 1.0300058, 1.0177208);
  4 = @(76.48186, 0.0, 0.0, 0.0, 0.0, 0.0);
  5 = @(8.503787599999999, 2.2529976, 1.1123126, 1.0977526, 1.3402676, 1.5331876);
  6 = @(0.8544899999999999, 0.807625, 0.0, 0.0, 0.372645, 0.0);
  7 = @(20.7405214, 3.1247414, 3.3144764, 2.8981514, 4.456071400000001, 3.0442064);
  8 = @(0.6889980000000001, 0.311803, 0.111603, 0.138903, 0.120703, 0.4728730000000001);
  9 = @(0.3578433999999999, 0.3259933999999999, 0.2440934, 0.2750333999999999, 0.2509184, 0.2809483999999999);
  10 = @(0.3899102, 0.3739852000000001, 0.4381402, 0.2961802000000001, 0.3052802000000001, 0.3580602);
  11 = @(0.1350142, 0.1409292, 0.0781392, 0.0672192, 0.0785942, 0.1832442);
  12 = @(0.107303026, 0.050473526, 0.029880226, 0.023592126, 0.032551076, 0.024019826);
  13 = @(0.103184908, 0.035890408, 0.029265608, 0.023282358, 0.06432790799999999, 0.024952208);
  14 = @(1432.1124, 1218.7174, 1201.4274, 1073.5724, 1634.1324, 1099.0524);
  15 = @(40831.37280000001, 35218.94780000002, 39697.05780000002, 33003.0978, 37431.61280000001, 33379.38280000001);
}

foreach ($row in $data.Keys) {
  $vals = $data[$row]
  for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item($row, $i + 2).Value = $vals[$i]
  }
}

Write-Output "Updated header row and data rows 2-15 across columns B:G"
